{"js": "// Update the worksheet date and the 25 division-problem answers.\n// Source duplicates text (e.g. \"71\u00f73=23, 2\" and \"29\u00f76=4, 5\" each occur\n// twice) map to DIFFERENT replacement values depending on where they sit\n// in the document, so replacements are done positionally (by paragraph /\n// table row+column) rather than by a global text search-and-replace.\n\nconst body = context.document.body;\n\n// --- 1. Update the title/date paragraph -------------------------------\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst dateParagraph = paragraphs.items[0];\ndateParagraph.getRange().insertText(\"2024-12-20 Friday\", \"Replace\");\n\n// --- 2. Update the five rows of division answers in the table ---------\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// Each inner array holds the five new cell values for one populated\n// table row, in the same top-to-bottom / left-to-right order as the\n// rows appear in the document (rows 0, 4, 8, 12, 16 of the 20-row table;\n// the rows between them are the blank answer rows and are left alone).\nconst newRowValues = [\n  [\"92\u00f76=15, 2\", \"18\u00f72=9, 0\", \"46\u00f79=5, 1\", \"26\u00f79=2, 8\", \"47\u00f75=9, 2\"],\n  [\"70\u00f77=10, 0\", \"97\u00f77=13, 6\", \"68\u00f78=8, 4\", \"73\u00f74=18, 1\", \"43\u00f75=8, 3\"],\n  [\"79\u00f78=9, 7\", \"44\u00f73=14, 2\", \"53\u00f77=7, 4\", \"75\u00f75=15, 0\", \"40\u00f72=20, 0\"],\n  [\"53\u00f74=13, 1\", \"30\u00f72=15, 0\", \"34\u00f72=17, 0\", \"39\u00f74=9, 3\", \"75\u00f73=25, 0\"],\n  [\"19\u00f73=6, 1\", \"74\u00f73=24, 2\", \"91\u00f78=11, 3\", \"11\u00f73=3, 2\", \"50\u00f73=16, 2\"],\n];\nconst rowIndexes = [0, 4, 8, 12, 16];\n\nfor (let r = 0; r < rowIndexes.length; r++) {\n  const tableRow = rowIndexes[r];\n  const values = newRowValues[r];\n  for (let c = 0; c < values.length; c++) {\n    table.getCell(tableRow, c).value = values[c];\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the worksheet date and the 25 division-problem answers.\n# Source duplicates text (e.g. \"71\u00f73=23, 2\" and \"29\u00f76=4, 5\" each occur\n# twice) map to DIFFERENT replacement values depending on where they sit\n# in the document, so replacements are done positionally (by paragraph /\n# table row+column) rather than by a global Find/Replace.\n\n$d = $word.ActiveDocument\n\n# --- 1. Update the title/date paragraph --------------------------------\n$d.Paragraphs(1).Range.Text = \"2024-12-20 Friday\"\n\n# --- 2. Update the five rows of division answers in the table ----------\n$t = $d.Tables(1)\n\n# Each inner array holds the five new cell values for one populated\n# table row, in the same top-to-bottom / left-to-right order as the rows\n# appear in the document (table rows 1, 5, 9, 13, 17 \u2014 1-based \u2014 of the\n# 20-row table; the rows in between are the blank answer rows and are\n# left untouched).\n$newRowValues = @(\n    ,@(\"92\u00f76=15, 2\", \"18\u00f72=9, 0\", \"46\u00f79=5, 1\", \"26\u00f79=2, 8\", \"47\u00f75=9, 2\")\n    ,@(\"70\u00f77=10, 0\", \"97\u00f77=13, 6\", \"68\u00f78=8, 4\", \"73\u00f74=18, 1\", \"43\u00f75=8, 3\")\n    ,@(\"79\u00f78=9, 7\", \"44\u00f73=14, 2\", \"53\u00f77=7, 4\", \"75\u00f75=15, 0\", \"40\u00f72=20, 0\")\n    ,@(\"53\u00f74=13, 1\", \"30\u00f72=15, 0\", \"34\u00f72=17, 0\", \"39\u00f74=9, 3\", \"75\u00f73=25, 0\")\n    ,@(\"19\u00f73=6, 1\", \"74\u00f73=24, 2\", \"91\u00f78=11, 3\", \"11\u00f73=3, 2\", \"50\u00f73=16, 2\")\n)\n$rowIndexes = @(1, 5, 9, 13, 17)\n\nfor ($r = 0; $r -lt $rowIndexes.Count; $r++) {\n    $tableRow = $rowIndexes[$r]\n    $values = $newRowValues[$r]\n    for ($c = 0; $c -lt $values.Count; $c++) {\n        $t.Cell($tableRow, $c + 1).Range.Text = $values[$c]\n    }\n}\n"}
